# Rewrite the body of the document to add the "matrices" exercises
# (N2 and N3) in the middle of the existing "Ejercicio N1" / counting
# exercise, and renumber that original exercise to N9, per the commit
# "add matrices reforzamiento N2".
#
# The new body is expressed as Word-processing XML fragments (one per
# paragraph) and injected in a single InsertXML call against the whole
# document Content range — this lets us control run/proofErr/bookmark
# boundaries exactly, the same way Word itself records them.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$paragraphs = @(
  # Ejercicio N1 heading (split across two runs, both bold)
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ejercicio N</w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1</w:t></w:r></w:p>",

  # Matrix 3x3 exercise body text
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t>Crear una matriz de 3×3 con los números del 1 al 9. Mostrar por pantalla, tal como aparece en la matriz.</w:t></w:r></w:p>",

  # Ejercicio N2 heading
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ejercicio N2</w:t></w:r></w:p>",

  # Matrix 5xn exercise body text
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t>Crear una matriz de 5 filas y n columnas (se pide al usuario). Rellenarlo con números aleatorios entre 0 y 10.</w:t></w:r></w:p>",

  # Ejercicio N3 heading
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ejercicio N3</w:t></w:r></w:p>",

  # Two nxn matrices exercise body text (with spellcheck markers + the
  # _GoBack bookmark left by the last edit position)
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=`"preserve`">Crear dos matrices de </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t>nxn</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=`"preserve`"> y sumar sus valores, los resultados se deben almacenar en otra matriz. Los valores y la longitud, </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t>seran</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=`"preserve`"> insertados por el usuario. Mostrar las matrices originales y el resultado.</w:t></w:r>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>",

  # Ejercicio N9 heading (the original "Ejercicio N1" content, renumbered,
  # split across two runs, both bold)
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ejercicio N</w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>9</w:t></w:r></w:p>",

  # Original array-counting exercise text (unchanged wording, merged runs)
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/></w:pPr>" +
  "<w:proofErr w:type=`"gramStart`"/><w:r><w:t>Programa java</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> que guarda en un array 10 números enteros que se leen por teclado. A continuación, se recorre el array y calcula cuantos números son positivos, cuantos negativos y cuantos ceros.</w:t></w:r></w:p>",

  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/></w:pPr>" +
  "<w:r><w:t>contar el número de elementos positivos, negativos y ceros de un array de 10 elementos.</w:t></w:r></w:p>",

  # Trailing empty paragraph (unchanged)
  "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/></w:pPr></w:p>"
)

$xml = [string]::Join("", $paragraphs)

$d.Content.InsertXML($xml)

Write-Output "Rebuilt body with $($paragraphs.Length) paragraphs"
